$wb = $excel.ActiveWorkbook

# --- Sheet: Forecast Comparison (Seasonality Index, column L) ---
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")

$wsForecast.Range("L2").Value = 0.99
$wsForecast.Range("L3").Value = 1
$wsForecast.Range("L4").Value = 0.86
$wsForecast.Range("L5").Value = 0.86
$wsForecast.Range("L6").Value = 0.9
$wsForecast.Range("L7").Value = 1.01
$wsForecast.Range("L8").Value = 0.8100000000000001
$wsForecast.Range("L9").Value = 0.86
$wsForecast.Range("L10").Value = 0.84
$wsForecast.Range("L11").Value = 0.99
$wsForecast.Range("L12").Value = 1.17
$wsForecast.Range("L13").Value = 1.05
$wsForecast.Range("L14").Value = 0.9
$wsForecast.Range("L15").Value = 1.08
$wsForecast.Range("L16").Value = 0.9399999999999999
$wsForecast.Range("L17").Value = 0.82

# --- Sheet: Summary (text values in column B) ---
# These cells hold numeric-looking text ("10", "5", "2", ...), so force
# text formatting while assigning, then restore the default style so we
# don't leave a new number format applied to the cell itself.
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("B9").NumberFormat = "@"
$wsSummary.Range("B9").Value = "12"
$wsSummary.Range("B9").Style = "Normal"

$wsSummary.Range("B10").NumberFormat = "@"
$wsSummary.Range("B10").Value = "6"
$wsSummary.Range("B10").Style = "Normal"

$wsSummary.Range("B11").NumberFormat = "@"
$wsSummary.Range("B11").Value = "3"
$wsSummary.Range("B11").Style = "Normal"
